# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp (10:31 -> 11:48)
# - Update case counters for several countries with the newer snapshot
# - A handful of countries swapped ranking position (their row's data
#   got overtaken by the neighboring row), so the country label + stats
#   for those rows are re-assigned accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Agosto de 2020 a las 11:48"

# --- India (row 6) ---
$ws.Range("B6").Value = 2091549
$ws.Range("C6").Value = 4685
$ws.Range("D6").Value = 1429178
$ws.Range("E6").Value = 619750
$ws.Range("G6").Value = 43
$ws.Range("H6").Value = 42621

# --- Banglades (row 18) ---
$ws.Range("B18").Value = 255113
$ws.Range("C18").Value = 2611
$ws.Range("D18").Value = 146604
$ws.Range("E18").Value = 105144
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = 3365

# --- Filipinas (row 25) ---
$ws.Range("B25").Value = 126885
$ws.Range("C25").Value = 4131
$ws.Range("D25").Value = 67117
$ws.Range("E25").Value = 57559
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = 2209

# --- Indonesia (row 26) ---
$ws.Range("B26").Value = 123503
$ws.Range("C26").Value = 2277
$ws.Range("D26").Value = 79306
$ws.Range("E26").Value = 38539
$ws.Range("G26").Value = 65
$ws.Range("H26").Value = 5658

# --- Oman (row 35) ---
$ws.Range("B35").Value = 81357
$ws.Range("C35").Value = 290
$ws.Range("D35").Value = 73481
$ws.Range("E35").Value = 7367
$ws.Range("G35").Value = 7
$ws.Range("H35").Value = 509

# --- Polonia (row 49) ---
$ws.Range("B49").Value = 51167
$ws.Range("C49").Value = 843
$ws.Range("D49").Value = 36403
$ws.Range("E49").Value = 12964
$ws.Range("G49").Value = 13
$ws.Range("H49").Value = 1800

# --- Consejo Danes para los Refugiados (row 87) ---
$ws.Range("B87").Value = 9436
$ws.Range("C87").Value = 81
$ws.Range("D87").Value = 8275
$ws.Range("E87").Value = 943

# --- Malasia (row 88) ---
$ws.Range("B88").Value = 9070
$ws.Range("C88").Value = 7
$ws.Range("D88").Value = 8775
$ws.Range("E88").Value = 170

# --- Finlandia (row 94) ---
$ws.Range("B94").Value = 7568
$ws.Range("C94").Value = 14
$ws.Range("E94").Value = 257

# --- Hong Kong (row 111) ---
$ws.Range("B111").Value = 4008
$ws.Range("C111").Value = 69
$ws.Range("D111").Value = 2755
$ws.Range("E111").Value = 1206

# --- Rows 123/124: Eslovaquia overtakes Mali ---
$ws.Range("A123").Value = "Eslovaquia"
$ws.Range("B123").Value = 2566
$ws.Range("C123").Value = 43
$ws.Range("D123").Value = 1861
$ws.Range("E123").Value = 674
$ws.Range("H123").Value = 31

$ws.Range("A124").Value = "Mali"
$ws.Range("B124").Value = 2561
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 1956
$ws.Range("E124").Value = 480
$ws.Range("H124").Value = 125

# --- Rows 127/128/129: Lituania overtakes Mozambique & Surinam ---
$ws.Range("A127").Value = "Lituania"
$ws.Range("B127").Value = 2231
$ws.Range("C127").Value = 37
$ws.Range("D127").Value = 1668
$ws.Range("E127").Value = 482
$ws.Range("H127").Value = 81

$ws.Range("A128").Value = "Mozambique"
$ws.Range("B128").Value = 2213
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 827
$ws.Range("E128").Value = 1371
$ws.Range("H128").Value = 15

$ws.Range("A129").Value = "Surinam"
$ws.Range("B129").Value = 2203
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 1505
$ws.Range("E129").Value = 669
$ws.Range("H129").Value = 29

# --- Estonia (row 130), same rank but refreshed numbers ---
$ws.Range("B130").Value = 2147
$ws.Range("C130").Value = 14
$ws.Range("D130").Value = 1961
$ws.Range("E130").Value = 123

# --- Rows 202/203: Timor Oriental overtakes Santa Lucia (tie, values unchanged) ---
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"
